$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, D, F, G, H, I, J, O across rows 2-25
# (columns A, E, K, L, M, N are unchanged by this edit)
$columns = @("B", "C", "D", "F", "G", "H", "I", "J", "O")

$newValues = @{
    2 = @(1.222612698847513, 0.3167053047682771, 0.4143797890649239, 0.8364779490808232, 0.2982334656762688, 0.4750511364174841, 0.3006879014162118, 0.3462298092455995, 1.473643830561599)
    3 = @(1.073878573760169, 0.2766488744388482, 0.4016173592216035, 0.8386635446224417, 0.3014336508809308, 0.4809699869310222, 0.309069743780551, 0.3344768212709965, 1.49254009397842)
    4 = @(0.9822093549041711, 0.2519445951081991, 0.3939300330000606, 0.8406840990580875, 0.303766264702432, 0.4849201392950206, 0.3145365653768306, 0.3274643494823266, 1.505573910301109)
    5 = @(0.9447692203120255, 0.2418506508352039, 0.3908350144149608, 0.8416778334776893, 0.3048089176181961, 0.4866092254590697, 0.3168446627244541, 0.3246579472330211, 1.511244291836348)
    6 = @(0.938547312496155, 0.2401729664424295, 0.3903233674980555, 0.8418531219083718, 0.304987600885319, 0.4868944880871737, 0.3172327651619957, 0.3241950397753897, 1.512207508906371)
    7 = @(0.9817047613000227, 0.2518085719190424, 0.39388813987388, 0.8406968115924229, 0.3037799538840531, 0.4849425976580051, 0.3145673682623675, 0.327426293952314, 1.50564893075493)
    8 = @(1.171402576358446, 0.3029170201437807, 0.4099485217311383, 0.8370905218483742, 0.2992603572311268, 0.4770263306468934, 0.3035113389170938, 0.342135030054024, 1.479861597534324)
    9 = @(1.540551026138019, 0.4022426811935702, 0.4426172444995871, 0.8354172772014508, 0.2933314174469999, 0.4640131813874433, 0.2843820624768298, 0.3726004020783904, 1.440693541859119)
    10 = @(1.809914975232971, 0.4746377333627265, 0.4673286673365169, 0.8375005616469622, 0.2907871219881244, 0.4559879935132756, 0.271896249695418, 0.3959797135995586, 1.418926843147446)
    11 = @(1.932032091250221, 0.5074399961117138, 0.4787235117588295, 0.8391722904748491, 0.2900279917828215, 0.4526716461062321, 0.2665594605712673, 0.4068337783279929, 1.41055960998041)
    12 = @(1.978212188632654, 0.5198419056442845, 0.4830603583512527, 0.8399097929076049, 0.2897981953008966, 0.4514640039606093, 0.2645881314343539, 0.4109754668990888, 1.407612858633982)
    13 = @(1.968269326581776, 0.5171718169251562, 0.4821253716033596, 0.8397463075068714, 0.2898451149098165, 0.4517219464832181, 0.2650104823905162, 0.4100820789906976, 1.408237615371803)
    14 = @(1.935832634929852, 0.5084607060308599, 0.4790798696820389, 0.8392308698699225, 0.2900079282346795, 0.4525713263143629, 0.2663962826259434, 0.4071738862017042, 1.41031272888506)
    15 = @(1.915955933768942, 0.5031223231403601, 0.4772172533549508, 0.8389287626145148, 0.2901151781389828, 0.4530978744910925, 0.2672515921201359, 0.4053966362678949, 1.411612704432187)
    16 = @(1.801925675882728, 0.4724913304020788, 0.4665870585289156, 0.8374059087349224, 0.2908447813881878, 0.4562114617969613, 0.2722519474364073, 0.3952747792514941, 1.419504632300402)
    17 = @(1.731862692796085, 0.4536661747603148, 0.4601049450990899, 0.8366573718863037, 0.2913946622679973, 0.4582072585495212, 0.2754075598417494, 0.3891213912392288, 1.424739873327866)
    18 = @(1.691525155972613, 0.4428261977390662, 0.4563910647746923, 0.8362949624410092, 0.2917484038354132, 0.4593866523320997, 0.2772548602125307, 0.3856027077237343, 1.427895420563942)
    19 = @(1.677860922143509, 0.4391538895376357, 0.4551360984943926, 0.8361839476829545, 0.2918745971061227, 0.4597913755455778, 0.2778858596885829, 0.3844148757481918, 1.428988598446026)
    20 = @(1.739325091253306, 0.4556714179889241, 0.4607934820925834, 0.8367300013288528, 0.2913322461775394, 0.4579915453634129, 0.2750682971983416, 0.3897742993861328, 1.424167622919555)
    21 = @(1.945361810524389, 0.5110199087277465, 0.4799738156135618, 0.8393794288200951, 0.2899585377222138, 0.4523205342424035, 0.2659878913509157, 0.4080272379676018, 1.409697191311764)
    22 = @(2.079650149088764, 0.5470787570482685, 0.4926366308290824, 0.8417199828360964, 0.2893970655337768, 0.448895118889169, 0.2603425332772238, 0.4201401541557317, 1.401532882039987)
    23 = @(2.008012614282791, 0.5278442231141298, 0.4858666626041099, 0.8404149497763882, 0.2896658293066494, 0.4506975880897173, 0.2633290170336258, 0.4136584530219096, 1.405771669184517)
    24 = @(1.735951519468301, 0.4547649008679855, 0.4604821547160611, 0.836696953971817, 0.2913603473767239, 0.4580889696743711, 0.2752215747644655, 0.3894790604950771, 1.424425883444897)
    25 = @(1.441003315566377, 0.3754719321167386, 0.4336544106836868, 0.8352897452290335, 0.2946188591388221, 0.4672642825949467, 0.289282535426385, 0.364184298307535, 1.450062801106597)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $vals[$i]
    }
}